# Apply updated cryptocurrency Price (D) and Volume(1h) (E) values
# exactly as captured by the upstream data refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep the exact text representation (no numeric
    # auto-conversion / trailing-zero loss), then restore the cell's
    # original (default) formatting so no stray style is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '27.763.19'
$ws.Range("E2").Value = '  -0.37%  '
Set-TextValue $ws.Range("D3") '1.615.78'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  +0.24%  '
Set-TextValue $ws.Range("D5") '210.74'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("E7").Value = '  +0.27%  '
Set-TextValue $ws.Range("D8") '22.89'
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("E9").Value = '  +0.07%  '
Set-TextValue $ws.Range("D10") '0.0602'
$ws.Range("E10").Value = '  -1.63%  '
Set-TextValue $ws.Range("D11") '0.0879'
$ws.Range("E11").Value = '  +0.07%  '
Set-TextValue $ws.Range("D12") '1.848.17'
$ws.Range("E12").Value = '  -0.42%  '
Set-TextValue $ws.Range("D13") '1.616.42'
$ws.Range("E13").Value = '  -0.48%  '
Set-TextValue $ws.Range("D14") '3.93'
$ws.Range("E14").Value = '  -2.24%  '
Set-TextValue $ws.Range("D15") '0.551'
$ws.Range("E15").Value = '  -1.87%  '
Set-TextValue $ws.Range("D16") '64.38'
$ws.Range("E16").Value = '  -1.53%  '
Set-TextValue $ws.Range("D17") '27.779.19'
$ws.Range("E17").Value = '  -0.31%  '
Set-TextValue $ws.Range("D18") '225.98'
$ws.Range("E18").Value = '  -1.65%  '
Set-TextValue $ws.Range("D19") '7.56'
$ws.Range("E19").Value = '  -1.14%  '
Set-TextValue $ws.Range("D20") '0.0₃0710'
$ws.Range("E20").Value = '  -1.63%  '
$ws.Range("E21").Value = '  +0.25%  '
Set-TextValue $ws.Range("D22") '4.29'
$ws.Range("E22").Value = '  -0.73%  '
Set-TextValue $ws.Range("D23") '9.95'
$ws.Range("E23").Value = '  -1.96%  '
Set-TextValue $ws.Range("D24") '2.05'
$ws.Range("E24").Value = '  +1.04%  '
Set-TextValue $ws.Range("D25") '154.68'
$ws.Range("E25").Value = '  -0.02%  '
Set-TextValue $ws.Range("D26") '6.88'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  -2.05%  '
Set-TextValue $ws.Range("D29") '15.27'
$ws.Range("E29").Value = '  -1.62%  '
Set-TextValue $ws.Range("D30") '1.17'
$ws.Range("E30").Value = '  -0.70%  '
Set-TextValue $ws.Range("D31") '0.0477'
$ws.Range("E31").Value = '  -0.72%  '
Set-TextValue $ws.Range("D32") '3.36'
$ws.Range("E32").Value = '  -1.41%  '
Set-TextValue $ws.Range("D33") '1.393.15'
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("E34").Value = '  -0.88%  '
Set-TextValue $ws.Range("D35") '1.59'
$ws.Range("E35").Value = '  +1.24%  '
Set-TextValue $ws.Range("D36") '0.972'
$ws.Range("E36").Value = '  -2.73%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  -1.29%  '
Set-TextValue $ws.Range("D39") '0.550'
$ws.Range("E39").Value = '  -0.87%  '
Set-TextValue $ws.Range("D40") '0.841'
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("E42").Value = '  -2.78%  '
Set-TextValue $ws.Range("D43") '65.04'
$ws.Range("E43").Value = '  -1.33%  '
Set-TextValue $ws.Range("D44") '5.34'
$ws.Range("E44").Value = '  -2.70%  '
Set-TextValue $ws.Range("D45") '1.76'
$ws.Range("E45").Value = '  -3.90%  '
Set-TextValue $ws.Range("D46") '1.756.84'
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("E47").Value = '  -3.24%  '
Set-TextValue $ws.Range("D48") '89.48'
$ws.Range("E49").Value = '  -2.37%  '
Set-TextValue $ws.Range("D50") '0.0990'
$ws.Range("E50").Value = '  -3.42%  '
$ws.Range("E51").Value = '  -0.48%  '
